# Add data for 2023-08-31
# Updates the 2023 (column J) crime-count figures across the citywide
# totals sheet, the "By Neighborhood" rollup sheet, and the individual
# neighborhood sheets, reflecting the refreshed extract through 2023-08-31.
# A handful of other-year cells (columns E/H/I) that were corrected in the
# same data refresh are also updated.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5042
$ws.Range("J3").Value = 5353
$ws.Range("E4").Value = 2003
$ws.Range("H4").Value = 1699
$ws.Range("I4").Value = 1772
$ws.Range("J4").Value = 1194
$ws.Range("J6").Value = 6667
$ws.Range("E7").Value = 26008
$ws.Range("H7").Value = 26010
$ws.Range("I7").Value = 26224
$ws.Range("J7").Value = 18676

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 48
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 249

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 211
$ws.Range("J3").Value = 298
$ws.Range("J7").Value = 808

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 97
$ws.Range("J3").Value = 95
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 109
$ws.Range("J7").Value = 287

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J3").Value = 30
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 178
$ws.Range("J3").Value = 274
$ws.Range("J4").Value = 55
$ws.Range("J7").Value = 719

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 149
$ws.Range("J4").Value = 67
$ws.Range("J5").Value = 59
$ws.Range("J7").Value = 549
$ws.Range("J8").Value = 1189
$ws.Range("J10").Value = 128
$ws.Range("J14").Value = 87
$ws.Range("J18").Value = 161
$ws.Range("J20").Value = 391
$ws.Range("J21").Value = 49
$ws.Range("J23").Value = 180
$ws.Range("J27").Value = 104
$ws.Range("J29").Value = 1043
$ws.Range("J30").Value = 75
$ws.Range("J33").Value = 854
$ws.Range("J37").Value = 587
$ws.Range("I42").Value = 1008
$ws.Range("J42").Value = 755
$ws.Range("J44").Value = 141
$ws.Range("J51").Value = 239
$ws.Range("J52").Value = 473
$ws.Range("J53").Value = 249
$ws.Range("J55").Value = 238
$ws.Range("E63").Value = 347
$ws.Range("J63").Value = 69
$ws.Range("J65").Value = 488
$ws.Range("J67").Value = 719
$ws.Range("J79").Value = 538
$ws.Range("J81").Value = 20
$ws.Range("J85").Value = 808
$ws.Range("H86").Value = 150
$ws.Range("J88").Value = 205
$ws.Range("J89").Value = 234
$ws.Range("J91").Value = 207
$ws.Range("J93").Value = 79
$ws.Range("J94").Value = 178
$ws.Range("J95").Value = 280
$ws.Range("J96").Value = 227
$ws.Range("J97").Value = 152
$ws.Range("J98").Value = 123
$ws.Range("J99").Value = 287
$ws.Range("E101").Value = 26008
$ws.Range("H101").Value = 26010
$ws.Range("I101").Value = 26224
$ws.Range("J101").Value = 18676

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 177
$ws.Range("J3").Value = 204
$ws.Range("J6").Value = 167
$ws.Range("J7").Value = 587

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 280
$ws.Range("J6").Value = 289
$ws.Range("J7").Value = 854

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 138
$ws.Range("J7").Value = 488

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 311
$ws.Range("J3").Value = 362
$ws.Range("J7").Value = 1043

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 45
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 141

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 175
$ws.Range("J3").Value = 165
$ws.Range("J4").Value = 22
$ws.Range("J7").Value = 549

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 167
$ws.Range("I4").Value = 57
$ws.Range("J6").Value = 386
$ws.Range("I7").Value = 1008
$ws.Range("J7").Value = 755

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 75
$ws.Range("J7").Value = 234

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 86
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 153
$ws.Range("J7").Value = 538

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 131
$ws.Range("J4").Value = 41
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 391

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 32
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 75
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 197
$ws.Range("J7").Value = 473

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 43
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J2").Value = 33
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("H4").Value = 69
$ws.Range("H7").Value = 150

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 56
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 227

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 333
$ws.Range("J3").Value = 359
$ws.Range("J4").Value = 72
$ws.Range("J6").Value = 391
$ws.Range("J7").Value = 1189

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 20

Write-Output "Updated 155 cells across 40 worksheets."
